# feat: add 2022-Q1 data
#
# Before: sheets = [ "2021-Q4", "总计" ]
# After:  sheets = [ "2021-Q4", "2022-Q1", "总计" ]
#   - "2021-Q4" is untouched.
#   - "2022-Q1" is a brand-new sheet (same layout as "2021-Q4") inserted
#     right after "2021-Q4" and before "总计".
#   - "总计" gets a new row 2 for "2022-Q1" (2 funds, 0.64 亿元) and the
#     existing "2021-Q4" row is pushed down to row 3 (its running index
#     becomes 1).

$wb = $excel.ActiveWorkbook

$xlCenter = -4108
$xlTop = -4160
$xlContinuous = 1

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q1" worksheet right after "2021-Q4".
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q1 = $wb.Worksheets.Add($null, $q4)
$q1.Name = "2022-Q1"

# Header row (bold, centered, top-aligned, thin box border) - mirrors
# the "2021-Q4" sheet's header styling.
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$hdr = $q1.Range("B1:H1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = $xlCenter
$hdr.VerticalAlignment = $xlTop
$hdr.Borders.LineStyle = $xlContinuous

# Row 2 - 004854 / 广发中证全指汽车指数A
$q1.Range("A2").Value = 0

$q1.Range("B2").NumberFormat = "@"
$q1.Range("B2").Value = "004854"
$q1.Range("C2").Value = "广发中证全指汽车指数A"
$q1.Range("D2").NumberFormat = "@"
$q1.Range("D2").Value = "22.01"
$q1.Range("E2").NumberFormat = "@"
$q1.Range("E2").Value = "94.43"
$q1.Range("F2").NumberFormat = "@"
$q1.Range("F2").Value = "2.27"
$q1.Range("G2").NumberFormat = "@"
$q1.Range("G2").Value = "0.4996"
$q1.Range("H2").Value = 10

# Row 3 - 004855 / 广发中证全指汽车指数C
$q1.Range("A3").Value = 1

$q1.Range("B3").NumberFormat = "@"
$q1.Range("B3").Value = "004855"
$q1.Range("C3").Value = "广发中证全指汽车指数C"
$q1.Range("D3").NumberFormat = "@"
$q1.Range("D3").Value = "6.11"
$q1.Range("E3").NumberFormat = "@"
$q1.Range("E3").Value = "94.43"
$q1.Range("F3").NumberFormat = "@"
$q1.Range("F3").Value = "2.27"
$q1.Range("G3").NumberFormat = "@"
$q1.Range("G3").Value = "0.1387"
$q1.Range("H3").Value = 10

# Index column (A2:A3) uses the same bold/border/centered style as the
# header row.
$idx = $q1.Range("A2:A3")
$idx.Font.Bold = $true
$idx.HorizontalAlignment = $xlCenter
$idx.VerticalAlignment = $xlTop
$idx.Borders.LineStyle = $xlContinuous

# ---------------------------------------------------------------------
# 2) Update the "总计" (totals) sheet: push the existing "2021-Q4" row
#    down to row 3 and write the new "2022-Q1" row into row 2.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Row 3 <- old row 2 ("2021-Q4", 2, 0.86), running index becomes 1.
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0.86

# Row 2 <- new "2022-Q1" summary row.
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.64

$totalIdx1 = $total.Range("A2")
$totalIdx1.Font.Bold = $true
$totalIdx1.HorizontalAlignment = $xlCenter
$totalIdx1.VerticalAlignment = $xlTop
$totalIdx1.Borders.LineStyle = $xlContinuous

$totalIdx2 = $total.Range("A3")
$totalIdx2.Font.Bold = $true
$totalIdx2.HorizontalAlignment = $xlCenter
$totalIdx2.VerticalAlignment = $xlTop
$totalIdx2.Borders.LineStyle = $xlContinuous

# ---------------------------------------------------------------------
# 3) Restore the original active sheet ("2021-Q4") / selection so the
#    workbook-level view state is unchanged by our edits.
# ---------------------------------------------------------------------
$q4.Activate()
[void]$q4.Range("A1").Select()
